$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("K2").Value = 1111
$ws.Range("K3").Value = 1056
$ws.Range("J4").Value = 274
$ws.Range("K4").Value = 230
$ws.Range("K5").Value = 64
$ws.Range("K6").Value = 1388
$ws.Range("J7").Value = 3927
$ws.Range("K7").Value = 3849

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("K2").Value = 66
$ws.Range("K4").Value = 14
$ws.Range("K7").Value = 227

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("K2").Value = 36
$ws.Range("K7").Value = 72

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("K2").Value = 43
$ws.Range("K6").Value = 47
$ws.Range("K7").Value = 157

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("K3").Value = 19
$ws.Range("K7").Value = 66

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("K2").Value = 26
$ws.Range("K3").Value = 40
$ws.Range("K5").Value = 6
$ws.Range("K6").Value = 45
$ws.Range("K7").Value = 125

$ws = $wb.Worksheets.Item('New City')
$ws.Range("K2").Value = 27
$ws.Range("K4").Value = 1
$ws.Range("K7").Value = 105

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("K3").Value = 26
$ws.Range("K7").Value = 72

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("K2").Value = 26
$ws.Range("K6").Value = 28
$ws.Range("K7").Value = 114
$ws.Range("K8").Value = 227
$ws.Range("K9").Value = 20
$ws.Range("K10").Value = 24
$ws.Range("K11").Value = 81
$ws.Range("K14").Value = 20
$ws.Range("K15").Value = 25
$ws.Range("K18").Value = 29
$ws.Range("K19").Value = 103
$ws.Range("K20").Value = 96
$ws.Range("K23").Value = 33
$ws.Range("K27").Value = 46
$ws.Range("K29").Value = 181
$ws.Range("K31").Value = 41
$ws.Range("K33").Value = 157
$ws.Range("K37").Value = 125
$ws.Range("K40").Value = 7
$ws.Range("K42").Value = 124
$ws.Range("K43").Value = 35
$ws.Range("K45").Value = 4
$ws.Range("K52").Value = 103
$ws.Range("K54").Value = 68
$ws.Range("K55").Value = 36
$ws.Range("K63").Value = 14
$ws.Range("K64").Value = 21
$ws.Range("K65").Value = 105
$ws.Range("K67").Value = 157
$ws.Range("J75").Value = 17
$ws.Range("K75").Value = 13
$ws.Range("K76").Value = 51
$ws.Range("K77").Value = 28
$ws.Range("K79").Value = 111
$ws.Range("K83").Value = 72
$ws.Range("K84").Value = 30
$ws.Range("K85").Value = 193
$ws.Range("K89").Value = 59
$ws.Range("K90").Value = 36
$ws.Range("K94").Value = 44
$ws.Range("K95").Value = 66
$ws.Range("K97").Value = 30
$ws.Range("K99").Value = 72
$ws.Range("J101").Value = 3927
$ws.Range("K101").Value = 3849

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("K2").Value = 20
$ws.Range("K7").Value = 41

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("K6").Value = 55
$ws.Range("K7").Value = 157

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("K6").Value = 8
$ws.Range("K7").Value = 30

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("K2").Value = 14
$ws.Range("K7").Value = 68

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("K2").Value = 47
$ws.Range("K3").Value = 55
$ws.Range("K7").Value = 181

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("K3").Value = 33
$ws.Range("K7").Value = 103

$ws = $wb.Worksheets.Item('River North')
$ws.Range("K2").Value = 9
$ws.Range("K4").Value = 5
$ws.Range("K6").Value = 27
$ws.Range("K7").Value = 51

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range("K2").Value = 7
$ws.Range("K3").Value = 3
$ws.Range("K7").Value = 20

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("K4").Value = 1
$ws.Range("K7").Value = 28

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("K2").Value = 29
$ws.Range("K7").Value = 124

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("K3").Value = 3
$ws.Range("K7").Value = 24

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("K2").Value = 15
$ws.Range("K7").Value = 36

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("K2").Value = 14
$ws.Range("K7").Value = 33

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("K2").Value = 41
$ws.Range("K4").Value = 9
$ws.Range("K6").Value = 25
$ws.Range("K7").Value = 111

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("K6").Value = 6
$ws.Range("K7").Value = 21

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("K3").Value = 29
$ws.Range("K7").Value = 96

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("K2").Value = 11
$ws.Range("K7").Value = 29

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("K2").Value = 44
$ws.Range("K6").Value = 32
$ws.Range("K7").Value = 114

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("K3").Value = 6
$ws.Range("K6").Value = 19
$ws.Range("K7").Value = 44

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("K2").Value = 7
$ws.Range("K7").Value = 25

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("K2").Value = 26
$ws.Range("K7").Value = 81

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range("K3").Value = 7
$ws.Range("K7").Value = 20

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("K6").Value = 8
$ws.Range("K7").Value = 26

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("K2").Value = 6
$ws.Range("K7").Value = 30

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("K3").Value = 23
$ws.Range("K7").Value = 59

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("K2").Value = 15
$ws.Range("K7").Value = 46

$ws = $wb.Worksheets.Item('Pullman')
$ws.Range("J4").Value = 2
$ws.Range("K6").Value = 2
$ws.Range("J7").Value = 17
$ws.Range("K7").Value = 13

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("K2").Value = 15
$ws.Range("K7").Value = 36

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("K6").Value = 17
$ws.Range("K7").Value = 35

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("K2").Value = 72
$ws.Range("K3").Value = 66
$ws.Range("K6").Value = 41
$ws.Range("K7").Value = 193

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("K2").Value = 14
$ws.Range("K7").Value = 28

$ws = $wb.Worksheets.Item('Jackson Park')
$ws.Range("K4").Value = 1
$ws.Range("K7").Value = 4

$ws = $wb.Worksheets.Item('Hegewisch')
$ws.Range("K2").Value = 3
$ws.Range("K7").Value = 7

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("K2").Value = 23
$ws.Range("K3").Value = 22
$ws.Range("K6").Value = 51
$ws.Range("K7").Value = 103
